$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 01:34"

# Estados Unidos (row 4) - updated case counts
$ws.Range("B4").Value = 1258051
$ws.Range("C4").Value = 20418
$ws.Range("D4").Value = 206203
$ws.Range("E4").Value = 977658
$ws.Range("F4").Value = 15808
$ws.Range("G4").Value = 1919
$ws.Range("H4").Value = 74190

# Colombia overtakes Republica Dominicana in the ranking, so rows 45/46 swap.
# Row 45 becomes Colombia with its updated figures.
$ws.Range("A45").Value = "Colombia"
$ws.Range("B45").Value = 8959
$ws.Range("C45").Value = 346
$ws.Range("D45").Value = 2148
$ws.Range("E45").Value = 6414
$ws.Range("F45").Value = 123
$ws.Range("G45").Value = 19
$ws.Range("H45").Value = 397

# Row 46 becomes Republica Dominicana with its updated figures.
$ws.Range("A46").Value = "Republica Dominicana"
$ws.Range("B46").Value = 8807
$ws.Range("C46").Value = 327
$ws.Range("D46").Value = 1905
$ws.Range("E46").Value = 6540
$ws.Range("F46").Value = 144
$ws.Range("G46").Value = 8
$ws.Range("H46").Value = 362
